$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 3290
$ws.Range("E2").Value = 306
$ws.Range("F2").Value = 306
$ws.Range("G2").Value = 142
$ws.Range("H2").Value = 115
$ws.Range("I2").Value = 115
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 13018
$ws.Range("L2").Value = 6542
$ws.Range("M2").Value = 6477
$ws.Range("N2").Value = 6477
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 137
$ws.Range("Q2").Value = 238
$ws.Range("R2").Value = -339
$ws.Range("S2").Value = 72
$ws.Range("T2").Value = 263
$ws.Range("U2").Value = -25
$ws.Range("V2").Value = 3739
$ws.Range("W2").Value = 9.3
$ws.Range("X2").Value = 3.49
$ws.Range("Y2").Value = 1.79
$ws.Range("Z2").Value = 0.88
$ws.Range("AA2").Value = 101
$ws.Range("AB2").Value = 4639.17
$ws.Range("AC2").Value = 419
$ws.Range("AD2").Value = 43.43
$ws.Range("AE2").Value = 23625
$ws.Range("AF2").Value = 0.77
$ws.Range("AG2").Value = 50
$ws.Range("AH2").Value = 0.27
$ws.Range("AI2").Value = 11.93
$ws.Range("AJ2").Value = 27415270

$ws.Range("D3").Value = 3576
$ws.Range("E3").Value = 390
$ws.Range("F3").Value = 390
$ws.Range("G3").Value = 233
$ws.Range("H3").Value = 167
$ws.Range("I3").Value = 167
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 13186
$ws.Range("L3").Value = 6557
$ws.Range("M3").Value = 6630
$ws.Range("N3").Value = 6630
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 137
$ws.Range("Q3").Value = 399
$ws.Range("R3").Value = -367
$ws.Range("S3").Value = -46
$ws.Range("T3").Value = 391
$ws.Range("U3").Value = 8
$ws.Range("V3").Value = 3735
$ws.Range("W3").Value = 10.89
$ws.Range("X3").Value = 4.67
$ws.Range("Y3").Value = 2.55
$ws.Range("Z3").Value = 1.27
$ws.Range("AA3").Value = 98.91
$ws.Range("AB3").Value = 4746.17
$ws.Range("AC3").Value = 609
$ws.Range("AD3").Value = 31.61
$ws.Range("AE3").Value = 24182
$ws.Range("AF3").Value = 0.8
$ws.Range("AG3").Value = 125
$ws.Range("AH3").Value = 0.65
$ws.Range("AI3").Value = 20.53
$ws.Range("AJ3").Value = 27415270

$ws.Range("D4").Value = 3774
$ws.Range("E4").Value = 434
$ws.Range("F4").Value = 434
$ws.Range("G4").Value = 324
$ws.Range("H4").Value = 295
$ws.Range("I4").Value = 295
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 13019
$ws.Range("L4").Value = 6221
$ws.Range("M4").Value = 6798
$ws.Range("N4").Value = 6799
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 137
$ws.Range("Q4").Value = 680
$ws.Range("R4").Value = -137
$ws.Range("S4").Value = -536
$ws.Range("T4").Value = 69
$ws.Range("U4").Value = 611
$ws.Range("V4").Value = 3268
$ws.Range("W4").Value = 11.51
$ws.Range("X4").Value = 7.81
$ws.Range("Y4").Value = 4.39
$ws.Range("Z4").Value = 2.25
$ws.Range("AA4").Value = 91.5
$ws.Range("AB4").Value = 4870.14
$ws.Range("AC4").Value = 1075
$ws.Range("AD4").Value = 14.23
$ws.Range("AE4").Value = 24799
$ws.Range("AF4").Value = 0.62
$ws.Range("AG4").Value = 180
$ws.Range("AH4").Value = 1.18
$ws.Range("AI4").Value = 16.74
$ws.Range("AJ4").Value = 27415270

$ws.Range("D5").Value = 3608
$ws.Range("E5").Value = 423
$ws.Range("F5").Value = 423
$ws.Range("G5").Value = 356
$ws.Range("H5").Value = 253
$ws.Range("I5").Value = 253
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 12740
$ws.Range("L5").Value = 5749
$ws.Range("M5").Value = 6991
$ws.Range("N5").Value = 6992
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 137
$ws.Range("Q5").Value = 577
$ws.Range("R5").Value = -196
$ws.Range("S5").Value = -361
$ws.Range("T5").Value = 372
$ws.Range("U5").Value = 205
$ws.Range("V5").Value = 2859
$ws.Range("W5").Value = 11.73
$ws.Range("X5").Value = 7.02
$ws.Range("Y5").Value = 3.67
$ws.Range("Z5").Value = 1.97
$ws.Range("AA5").Value = 82.23
$ws.Range("AB5").Value = 5018.82
$ws.Range("AC5").Value = 923
$ws.Range("AD5").Value = 15.05
$ws.Range("AE5").Value = 25502
$ws.Range("AF5").Value = 0.55
$ws.Range("AG5").Value = 180
$ws.Range("AH5").Value = 1.29
$ws.Range("AI5").Value = 19.49
$ws.Range("AJ5").Value = 27415270

$ws.Range("D6").Value = 3514
$ws.Range("E6").Value = 410
$ws.Range("F6").Value = 410
$ws.Range("G6").Value = 272
$ws.Range("H6").Value = 202
$ws.Range("I6").Value = 202
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 13018
$ws.Range("L6").Value = 5930
$ws.Range("M6").Value = 7088
$ws.Range("N6").Value = 7088
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 137
$ws.Range("Q6").Value = 121
$ws.Range("R6").Value = -257
$ws.Range("S6").Value = 103
$ws.Range("T6").Value = 152
$ws.Range("U6").Value = -31
$ws.Range("V6").Value = 3079
$ws.Range("W6").Value = 11.66
$ws.Range("X6").Value = 5.75
$ws.Range("Y6").Value = 2.87
$ws.Range("Z6").Value = 1.57
$ws.Range("AA6").Value = 83.67
$ws.Range("AB6").Value = 5130.13
$ws.Range("AC6").Value = 737
$ws.Range("AD6").Value = 14.66
$ws.Range("AE6").Value = 26282
$ws.Range("AF6").Value = 0.41
$ws.Range("AG6").Value = 125
$ws.Range("AH6").Value = 1.16
$ws.Range("AI6").Value = 16.7
$ws.Range("AJ6").Value = 27415270

$ws.Range("D7:AJ9").ClearContents()

Write-Host "Done"
